# Update the "Förändrad" (changed) date column (C) from 45179 to 45180
# for all data rows (rows 2 through 115) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSerial = 45179
$newSerial = 45180

$lastRow = 115
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value2 = $newSerial
    }
}
